$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct selection scopes: clear the G:H values that were incorrectly
# included for these rows (odd rows 3,5,7,9,11,13,15,17).
$rows = @(3, 5, 7, 9, 11, 13, 15, 17)
foreach ($r in $rows) {
    $ws.Range("G$r`:H$r").ClearContents()
}
